# Applies the "ARI_results" revision:
#  - updates the ARI score column (B2:B77) with re-run values
#  - applies a custom 4-decimal number format to the later block (B80:B120)
#  - adds a new blank separator row (121) styled like the existing one at row 78
#  - moves the saved selection to F80
#  - sets the print page setup (portrait / paper size 9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- B2:B77 updated values ----
$newValues = @{
    "B2"  = 0.3717
    "B3"  = 0.5125
    "B4"  = 0.3952
    "B5"  = 0.4994
    "B6"  = 0.5517
    "B8"  = 0.0323
    "B9"  = 0.0019
    "B10" = 0.9463
    "B11" = 1.0
    "B12" = 0.1921
    "B13" = 0.2888
    "B14" = 0.761
    "B15" = 0.4641
    "B16" = 0.0085
    "B17" = 0.2041
    "B18" = 0.0024
    "B19" = 0.0023
    "B20" = 0.0255
    "B21" = 0.2832
    "B22" = 0.1693
    "B23" = -0.0618
    "B24" = 0.3177
    "B25" = 0.003
    "B26" = 0.5501
    "B27" = 0.5549
    "B28" = 0.6162
    "B29" = 0.9814
    "B30" = 0.8853
    "B31" = 0.7948
    "B32" = 0.2472
    "B33" = 0.5247
    "B34" = 0.5376
    "B35" = 0.0521
    "B36" = 0.015
    "B37" = 0.0731
    "B38" = 0.0892
    "B39" = 0.0194
    "B40" = 0.4359
    "B41" = 0.4033
    "B42" = 0.8753
    "B43" = 0.4093
    "B44" = -0.0054
    "B45" = 0.295
    "B46" = 0.5742
    "B47" = 0.7253
    "B48" = 0.5457
    "B49" = 0.6386
    "B50" = 0.6312
    "B51" = -0.0005
    "B52" = 0.0641
    "B53" = 0.6723
    "B54" = 0.3937
    "B55" = 0.2319
    "B56" = 0.5717
    "B57" = 0.0383
    "B58" = 0.124
    "B59" = 0.0336
    "B60" = 0.0775
    "B61" = 0.2228
    "B62" = 0.5209
    "B63" = 0.1573
    "B64" = 0.2265
    "B65" = 0.8452
    "B66" = 0.3739
    "B67" = 0.543
    "B68" = 0.9612
    "B69" = 0.0928
    "B70" = 0.9897
    "B71" = 0.9968
    "B72" = 0.4787
    "B73" = 0.7206
    "B74" = 0.5905
    "B75" = 0.4707
    "B76" = 0.0098
    "B77" = 0.0877
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# ---- apply the new 4-decimal number format to the later results block ----
$ws.Range("B80:B120").NumberFormat = "0.0000_ "

# ---- new blank separator row 121, styled like the existing separator (row 78) ----
$ws.Range("B78").Copy()
$ws.Range("B121").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- move the remembered selection ----
[void]$ws.Range("F80").Select()

# ---- page setup for printing ----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "edit applied"
